$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 28.401376
$ws.Range("H2").Value = 85.204128
$ws.Range("I2").Value = 0.2813463917610605
$ws.Range("J2").Value = 0.2813463917610605
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 91.65675166666666
$ws.Range("N2").Value = 274.970255
$ws.Range("O2").Value = 0.9732857812565366
$ws.Range("P2").Value = 0.9732857812565366
$ws.Range("Q2").Value = 2603.177867023627
$ws.Range("R2").Value = 23428.60080321264
$ws.Range("S2").Value = 0.2738304427088714
$ws.Range("T2").Value = 0.2738304427088714

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 28.401376
$ws.Range("H3").Value = 85.204128
$ws.Range("I3").Value = 0.2813463917610605
$ws.Range("J3").Value = 0.2813463917610605
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6377586666666667
$ws.Range("N3").Value = 1.913276
$ws.Range("O3").Value = 0.006772239151537979
$ws.Range("P3").Value = 0.006772239151537978
$ws.Range("Q3").Value = 18.11322368925867
$ws.Range("R3").Value = 163.019013203328
$ws.Range("S3").Value = 0.001905345049428196
$ws.Range("T3").Value = 0.001905345049428196

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 28.401376
$ws.Range("H4").Value = 85.204128
$ws.Range("I4").Value = 0.2813463917610605
$ws.Range("J4").Value = 0.2813463917610605
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.877986
$ws.Range("N4").Value = 5.633958
$ws.Range("O4").Value = 0.01994197959192537
$ws.Range("P4").Value = 0.01994197959192537
$ws.Range("Q4").Value = 53.337386508736
$ws.Range("R4").Value = 480.036478578624
$ws.Range("S4").Value = 0.005610604002760908
$ws.Range("T4").Value = 0.005610604002760909

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.83037466666667
$ws.Range("H5").Value = 44.491124
$ws.Range("I5").Value = 0.1469109243485705
$ws.Range("J5").Value = 0.1469109243485705
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 91.65675166666666
$ws.Range("N5").Value = 274.970255
$ws.Range("O5").Value = 0.9732857812565366
$ws.Range("P5").Value = 0.9732857812565366
$ws.Range("Q5").Value = 1359.303967946291
$ws.Range("R5").Value = 12233.73571151662
$ws.Range("S5").Value = 0.1429863137797184
$ws.Range("T5").Value = 0.1429863137797184

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.83037466666667
$ws.Range("H6").Value = 44.491124
$ws.Range("I6").Value = 0.1469109243485705
$ws.Range("J6").Value = 0.1469109243485705
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6377586666666667
$ws.Range("N6").Value = 1.913276
$ws.Range("O6").Value = 0.006772239151537979
$ws.Range("P6").Value = 0.006772239151537978
$ws.Range("Q6").Value = 9.458199973580445
$ws.Range("R6").Value = 85.123799762224
$ws.Range("S6").Value = 0.0009949159136620235
$ws.Range("T6").Value = 0.0009949159136620233

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 14.83037466666667
$ws.Range("H7").Value = 44.491124
$ws.Range("I7").Value = 0.1469109243485705
$ws.Range("J7").Value = 0.1469109243485705
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.877986
$ws.Range("N7").Value = 5.633958
$ws.Range("O7").Value = 0.01994197959192537
$ws.Range("P7").Value = 0.01994197959192537
$ws.Range("Q7").Value = 27.85123599875466
$ws.Range("R7").Value = 250.661123988792
$ws.Range("S7").Value = 0.002929694655190085
$ws.Range("T7").Value = 0.002929694655190085

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 57.71632199999999
$ws.Range("H8").Value = 173.148966
$ws.Range("I8").Value = 0.571742683890369
$ws.Range("J8").Value = 0.571742683890369
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 91.65675166666666
$ws.Range("N8").Value = 274.970255
$ws.Range("O8").Value = 0.9732857812565366
$ws.Range("P8").Value = 0.9732857812565366
$ws.Range("Q8").Value = 5290.090592667369
$ws.Range("R8").Value = 47610.81533400632
$ws.Range("S8").Value = 0.5564690247679468
$ws.Range("T8").Value = 0.5564690247679468

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 57.71632199999999
$ws.Range("H9").Value = 173.148966
$ws.Range("I9").Value = 0.571742683890369
$ws.Range("J9").Value = 0.571742683890369
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6377586666666667
$ws.Range("N9").Value = 1.913276
$ws.Range("O9").Value = 0.006772239151537979
$ws.Range("P9").Value = 0.006772239151537978
$ws.Range("Q9").Value = 36.809084563624
$ws.Range("R9").Value = 331.281761072616
$ws.Range("S9").Value = 0.003871978188447759
$ws.Range("T9").Value = 0.003871978188447759

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 57.71632199999999
$ws.Range("H10").Value = 173.148966
$ws.Range("I10").Value = 0.571742683890369
$ws.Range("J10").Value = 0.571742683890369
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.877986
$ws.Range("N10").Value = 5.633958
$ws.Range("O10").Value = 0.01994197959192537
$ws.Range("P10").Value = 0.01994197959192537
$ws.Range("Q10").Value = 108.390444687492
$ws.Range("R10").Value = 975.5140021874278
$ws.Range("S10").Value = 0.01140168093397437
$ws.Range("T10").Value = 0.01140168093397437

